# Chapter 8 finished - add CART / Gini impurity bullets under "Decision Trees",
# keeping "Ensemble Learning and Random Forests" as the final bullet.
#
# Before:
#   ... Decision Trees (lvl0)
#       Ensemble Learning and Random Forests (lvl0)   <- holds the _GoBack bookmark
#
# After:
#   ... Decision Trees (lvl0)
#       Classification and Regression Tree (CART) algorithm (lvl1)
#       Gini impurity vs entropy (lvl1)                <- now holds the _GoBack bookmark
#       Ensemble Learning and Random Forests (lvl0)

$d = $word.ActiveDocument

# The last paragraph in the document is "Ensemble Learning and Random Forests";
# it carries the hidden _GoBack bookmark left over from the last edit.
$origPara = $d.Paragraphs.Last

# 1) Add an empty paragraph after it that will eventually become the new, final
#    "Ensemble Learning and Random Forests" bullet. Seed it with a sentinel
#    character so it is non-empty (merging into a non-empty paragraph later
#    keeps the run list clean).
$r = $origPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$spacer = $d.Paragraphs.Item($n)
$spacer.Range.InsertBefore([char]1)

# 2) Turn the original (bookmarked) paragraph's text into "Gini impurity vs
#    entropy" - the bookmark tags sit after the run and are untouched by the
#    text substitution.
$null = $origPara.Range.Find.Execute("Ensemble Learning and Random Forests", $true, $false, $false, $false, $false, $true, 1, $false, "Gini impurity vs entropy", 2)

# 3) Delete the (now Gini) paragraph's own paragraph mark so its content
#    merges forward into the sentinel paragraph. The surviving paragraph mark
#    is the sentinel's (freshly created, so it carries no rsid attributes),
#    while the bookmark and text are preserved.
$markRange = $d.Range($origPara.Range.End - 1, $origPara.Range.End)
$markRange.Delete()
$mergedIndex = $n - 1
$merged = $d.Paragraphs.Item($mergedIndex)

# 4) Strip the sentinel character back out.
$null = $merged.Range.Find.Execute([char]1, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 5) Demote the merged "Gini impurity vs entropy" paragraph to list level 1.
$merged.Range.ListFormat.ListIndent()

# 6) Insert the new "Classification and Regression Tree (CART) algorithm"
#    bullet right before the "Gini impurity vs entropy" paragraph. It
#    inherits the Gini paragraph's (already level-1) formatting, so no
#    further demotion is needed.
$insertPoint = $d.Range($merged.Range.Start, $merged.Range.Start)
$insertPoint.InsertBefore("Classification and Regression Tree (CART) algorithm" + [char]13)

# 7) Insert the final "Ensemble Learning and Random Forests" bullet after the
#    "Gini impurity vs entropy" paragraph, then promote it back to list
#    level 0.
$giniPara = $d.Paragraphs.Item($mergedIndex + 1)
$gr = $giniPara.Range
$gr.Collapse(0)
$gr.InsertParagraphAfter()
$ensemblePara = $d.Paragraphs.Item($mergedIndex + 2)
$ensemblePara.Range.InsertBefore("Ensemble Learning and Random Forests")
$ensemblePara.Range.ListFormat.ListOutdent()
